$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.852.00'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '3.438.27'
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.22'
$ws.Range("E5").Value = '  -1.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.25'
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").Value = '3.438.26'
$ws.Range("E7").Value = '  -0.16%  '

$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.75'
$ws.Range("E10").Value = '  +1.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.124'
$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.403'
$ws.Range("E12").Value = '  +2.33%  '

$ws.Range("D13").Value = '4.025.73'
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("E14").Value = '  +2.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.88'
$ws.Range("E15").Value = '  -2.02%  '

$ws.Range("D16").Value = '3.438.54'
$ws.Range("E16").Value = '  -0.23%  '

$ws.Range("E17").Value = '  -1.04%  '

$ws.Range("D18").Value = '62.889.26'
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("E19").Value = '  +1.69%  '

$ws.Range("E20").Value = '  +0.36%  '

$ws.Range("E21").Value = '  -1.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.42'
$ws.Range("E22").Value = '  -2.53%  '

$ws.Range("E23").Value = '  -0.62%  '

$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("E25").Value = '  -0.16%  '

$ws.Range("D26").Value = '3.583.39'
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("E27").Value = '  -3.73%  '

$ws.Range("E28").Value = '  -5.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.55'
$ws.Range("E29").Value = '  -2.01%  '

$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.06'
$ws.Range("E31").Value = '  -1.48%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.11'
$ws.Range("E32").Value = '  -2.07%  '

$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.25'
$ws.Range("E34").Value = '  -2.09%  '

$ws.Range("E35").Value = '  -8.93%  '

$ws.Range("E36").Value = '  -1.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.06'
$ws.Range("E37").Value = '  -0.17%  '

$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.62'
$ws.Range("E38").Value = '  +3.11%  '

$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.58'
$ws.Range("E39").Value = '  -0.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '168.74'
$ws.Range("E40").Value = '  +0.41%  '

$ws.Range("D41").Value = '3.475.47'
$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("E43").Value = '  -0.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.33'
$ws.Range("E44").Value = '  -1.26%  '

$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("E46").Value = '  -0.96%  '

$ws.Range("E47").Value = '  -3.26%  '

$ws.Range("D48").Value = '2.562.74'
$ws.Range("E48").Value = '  +1.77%  '

$ws.Range("E49").Value = '  +4.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.80'
$ws.Range("E50").Value = '  +0.99%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.63'
$ws.Range("E51").Value = '  -4.11%  '

